$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 2.52
$ws.Range("I2").Value = 2.62
$ws.Range("Q2").Value = 1.65
$ws.Range("V2").Value = 1.61
$ws.Range("W2").Value = 1.53
$ws.Range("AB2").Value = 16
$ws.Range("AE2").Value = 29
$ws.Range("F3").Value = 1.97
$ws.Range("I3").Value = 3.8
$ws.Range("K3").Value = 8.6
$ws.Range("L3").Value = 1.2
$ws.Range("V3").Value = 1.36
$ws.Range("K5").Value = 950
$ws.Range("P5").Value = 1.25
$ws.Range("Q5").Value = 1.45
$ws.Range("R5").Value = 1.18
$ws.Range("S5").Value = 1.46
$ws.Range("F6").Value = 3.6
$ws.Range("G6").Value = 3.95
$ws.Range("H6").Value = 2.08
$ws.Range("I6").Value = 2.24
$ws.Range("J6").Value = 3.45
$ws.Range("L6").Value = 1.4
$ws.Range("P6").Value = 1.96
$ws.Range("R6").Value = 1.37
$ws.Range("T6").Value = 1.75
$ws.Range("U6").Value = 2.14
$ws.Range("V6").Value = 1.81
$ws.Range("W6").Value = 1.34
$ws.Range("Y6").Value = 11.5
$ws.Range("Z6").Value = 17
$ws.Range("AA6").Value = 27
$ws.Range("AB6").Value = 15.5
$ws.Range("AD6").Value = 12.5
$ws.Range("AI6").Value = 44
$ws.Range("AJ6").Value = 75
$ws.Range("AN6").Value = 55
$ws.Range("AO6").Value = 19
$ws.Range("D8").Value = "Landskrona"
$ws.Range("E8").Value = "Varbergs BoIS"
$ws.Range("F8").Value = 2.76
$ws.Range("G8").Value = 3.05
$ws.Range("H8").Value = 2.48
$ws.Range("I8").Value = 2.7
$ws.Range("J8").Value = 3.7
$ws.Range("K8").Value = 3.8
$ws.Range("L8").Value = 1.39
$ws.Range("M8").Value = 1.07
$ws.Range("N8").Value = 2.08
$ws.Range("O8").Value = 1.08
$ws.Range("P8").Value = 1.3
$ws.Range("Q8").Value = 1.01
$ws.Range("R8").Value = 1.21
$ws.Range("S8").Value = 1.05
$ws.Range("T8").Value = 1.61
$ws.Range("U8").Value = 2.12
$ws.Range("V8").Value = 1.58
$ws.Range("W8").Value = 1.5
$ws.Range("X8").Value = 16
$ws.Range("Y8").Value = 13.5
$ws.Range("Z8").Value = 18.5
$ws.Range("AA8").Value = 40
$ws.Range("AB8").Value = 1000
$ws.Range("AC8").Value = 8.4
$ws.Range("AD8").Value = 12.5
$ws.Range("AE8").Value = 27
$ws.Range("AF8").Value = 21
$ws.Range("AG8").Value = 12.5
$ws.Range("AH8").Value = 16.5
$ws.Range("AI8").Value = 40
$ws.Range("AJ8").Value = 1000
$ws.Range("AK8").Value = 32
$ws.Range("AL8").Value = 44
$ws.Range("AM8").Value = 80
$ws.Range("AN8").Value = 25
$ws.Range("AO8").Value = 22
$ws.Range("A9").Value = "Swedish Superettan"
$ws.Range("C9").Value = "14:00:00"
$ws.Range("D9").Value = "Sandvikens"
$ws.Range("E9").Value = "Orgryte"
$ws.Range("F9").Value = 3.4
$ws.Range("G9").Value = 3.7
$ws.Range("H9").Value = 2.02
$ws.Range("I9").Value = 2.14
$ws.Range("J9").Value = 4
$ws.Range("K9").Value = 4.5
$ws.Range("L9").Value = 1.23
$ws.Range("M9").Value = 1.03
$ws.Range("N9").Value = 5.5
$ws.Range("O9").Value = 1.19
$ws.Range("P9").Value = 2.56
$ws.Range("Q9").Value = 1.55
$ws.Range("R9").Value = 1.62
$ws.Range("S9").Value = 2.36
$ws.Range("T9").Value = 1.55
$ws.Range("U9").Value = 2.58
$ws.Range("V9").Value = 1.87
$ws.Range("W9").Value = 1.37
$ws.Range("X9").Value = 32
$ws.Range("Y9").Value = 17
$ws.Range("Z9").Value = 18
$ws.Range("AA9").Value = 26
$ws.Range("AB9").Value = 24
$ws.Range("AC9").Value = 12
$ws.Range("AD9").Value = 13.5
$ws.Range("AE9").Value = 21
$ws.Range("AF9").Value = 30
$ws.Range("AG9").Value = 21
$ws.Range("AH9").Value = 18
$ws.Range("AI9").Value = 32
$ws.Range("AJ9").Value = 75
$ws.Range("AK9").Value = 42
$ws.Range("AL9").Value = 44
$ws.Range("AM9").Value = 60
$ws.Range("AN9").Value = 29
$ws.Range("AO9").Value = 11.5
$ws.Range("A10").Value = "Romanian Liga I"
$ws.Range("C10").Value = "14:30:00"
$ws.Range("D10").Value = "Botosani"
$ws.Range("E10").Value = "UTA Arad"
$ws.Range("F10").Value = 2.4
$ws.Range("G10").Value = 2.58
$ws.Range("H10").Value = 3.2
$ws.Range("I10").Value = 3.55
$ws.Range("J10").Value = 3.25
$ws.Range("K10").Value = 3.35
$ws.Range("L10").Value = 1.35
$ws.Range("N10").Value = 3.4
$ws.Range("O10").Value = 1.34
$ws.Range("P10").Value = 1.82
$ws.Range("Q10").Value = 1.99
$ws.Range("R10").Value = 1.32
$ws.Range("S10").Value = 3.55
$ws.Range("T10").Value = 1.76
$ws.Range("U10").Value = 2.06
$ws.Range("V10").Value = 1.4
$ws.Range("W10").Value = 1.64
$ws.Range("X10").Value = 16
$ws.Range("Y10").Value = 15
$ws.Range("Z10").Value = 28
$ws.Range("AA10").Value = 70
$ws.Range("AB10").Value = 12
$ws.Range("AC10").Value = 9.4
$ws.Range("AD10").Value = 17
$ws.Range("AE10").Value = 48
$ws.Range("AF10").Value = 19
$ws.Range("AG10").Value = 14
$ws.Range("AH10").Value = 22
$ws.Range("AI10").Value = 65
$ws.Range("AJ10").Value = 42
$ws.Range("AK10").Value = 34
$ws.Range("AL10").Value = 50
$ws.Range("AM10").Value = 120
$ws.Range("AN10").Value = 27
$ws.Range("AO10").Value = 48
$ws.Range("F11").Value = 2.94
$ws.Range("G11").Value = 3.05
$ws.Range("H11").Value = 2.2
$ws.Range("I11").Value = 2.36
$ws.Range("J11").Value = 4
$ws.Range("K11").Value = 4.8
$ws.Range("O11").Value = 1.19
$ws.Range("P11").Value = 2.44
$ws.Range("T11").Value = 1.53
$ws.Range("V11").Value = 1.73
$ws.Range("W11").Value = 1.48
$ws.Range("K12").Value = 950
$ws.Range("P12").Value = 1.24
$ws.Range("W12").Value = 3.35
$ws.Range("F13").Value = 2.78
$ws.Range("H13").Value = 2.72
$ws.Range("I13").Value = 2.86
$ws.Range("K13").Value = 3.45
$ws.Range("V13").Value = 1.53
$ws.Range("F14").Value = 1.48
$ws.Range("L14").Value = 1.37
$ws.Range("T14").Value = 2.02
$ws.Range("U14").Value = 1.68
$ws.Range("F15").Value = 3.4
$ws.Range("G15").Value = 3.7
$ws.Range("H15").Value = 2.52
$ws.Range("I15").Value = 2.68
$ws.Range("J15").Value = 2.92
$ws.Range("K15").Value = 3.1
$ws.Range("M15").Value = 1.18
$ws.Range("N15").Value = 2.14
$ws.Range("O15").Value = 1.77
$ws.Range("P15").Value = 1.38
$ws.Range("S15").Value = 8
$ws.Range("U15").Value = 1.57
$ws.Range("V15").Value = 1.59
$ws.Range("W15").Value = 1.37
$ws.Range("X15").Value = 6.8
$ws.Range("Y15").Value = 6.6
$ws.Range("Z15").Value = 14
$ws.Range("AA15").Value = 90
$ws.Range("AB15").Value = 8
$ws.Range("AD15").Value = 19
$ws.Range("AE15").Value = 160
$ws.Range("AF15").Value = 29
$ws.Range("AG15").Value = 24
$ws.Range("AH15").Value = 80
$ws.Range("AI15").Value = 1000
$ws.Range("AL15").Value = 1000
$ws.Range("AM15").Value = 430
$ws.Range("AN15").Value = 1000
$ws.Range("AO15").Value = 390
$ws.Range("F16").Value = 1.99
$ws.Range("G16").Value = 2.14
$ws.Range("H16").Value = 4.1
$ws.Range("I16").Value = 5.3
$ws.Range("K16").Value = 3.9
$ws.Range("N16").Value = 2.48
$ws.Range("O16").Value = 1.51
$ws.Range("Q16").Value = 2.46
$ws.Range("S16").Value = 4.9
$ws.Range("T16").Value = 2.02
$ws.Range("V16").Value = 1.23
$ws.Range("W16").Value = 1.87
$ws.Range("Y16").Value = 17
$ws.Range("AB16").Value = 8.800000000000001
$ws.Range("F18").Value = 1.75
$ws.Range("G18").Value = 1.9
$ws.Range("P18").Value = 1.7
$ws.Range("Q18").Value = 2.16
$ws.Range("R18").Value = 1.25
$ws.Range("T18").Value = 2.06
$ws.Range("U18").Value = 1.8
$ws.Range("AB18").Value = 8.6
$ws.Range("AH18").Value = 26
$ws.Range("F19").Value = 3.25
$ws.Range("H19").Value = 2.38
$ws.Range("I19").Value = 2.68
$ws.Range("K19").Value = 3.7
$ws.Range("R19").Value = 1.25
$ws.Range("T19").Value = 1.92
$ws.Range("V19").Value = 1.62
$ws.Range("AJ19").Value = 85
$ws.Range("AN19").Value = 70

Write-Output "Applied 239 cell updates"
